# Update the "Contact person(s)/Focal point" block (rows 7-10) with the
# new reporter's details, replacing the previous contact's information.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "Kalymbetova Yryskan"
$ws.Range("B8").Value = "yryskan.kalymbetova@gmail.com "
$ws.Range("B9").Value = "(0312) 32 46 55"
$ws.Range("B10").Value = "www.stat.gov.kg"

# Move/leave the active selection on the updated cell.
$ws.Range("B7").Select()
